$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text formatting (values like "30.375.12" must not be
# auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.375.12"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "1.937.08"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "0.7737"
$ws.Range("E5").Value = "  +7.09%  "

$ws.Range("D6").Value = "248.90"

$ws.Range("D7").Value = "0.9990"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").Value = "28.25"
$ws.Range("E8").Value = "  +1.59%  "

$ws.Range("D9").Value = "0.3221"
$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").Value = "0.07112"
$ws.Range("E10").Value = "  -2.15%  "

$ws.Range("D11").Value = "0.7902"
$ws.Range("E11").Value = "  -2.56%  "

$ws.Range("D12").Value = "0.08013"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("D13").Value = "1.936.66"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").Value = "5.395"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").Value = "95.00"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "14.66"
$ws.Range("E16").Value = "  -3.15%  "

$ws.Range("D17").Value = "30.391.03"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").Value = "256.82"
$ws.Range("E18").Value = "  +2.01%  "

$ws.Range("D19").Value = "0.000008043"
$ws.Range("E19").Value = "  -3.02%  "

$ws.Range("D20").Value = "5.819"
$ws.Range("E20").Value = "  -1.23%  "

$ws.Range("D21").Value = "2.191.98"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "0.9987"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").Value = "0.9982"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "6.827"
$ws.Range("E24").Value = "  -2.55%  "

$ws.Range("D25").Value = "9.649"

$ws.Range("D26").Value = "164.64"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.1357"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "19.20"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").Value = "2.317"
$ws.Range("E29").Value = "  -2.90%  "

$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("D31").Value = "1.532"
$ws.Range("E31").Value = "  -2.26%  "

$ws.Range("D32").Value = "4.453"
$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("D33").Value = "4.170"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").Value = "0.05226"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("D35").Value = "1.296"
$ws.Range("E35").Value = "  +0.59%  "

$ws.Range("D36").Value = "0.7552"
$ws.Range("E36").Value = "  +0.61%  "

$ws.Range("D37").Value = "2.767"
$ws.Range("E37").Value = "  +0.91%  "

$ws.Range("D38").Value = "0.01979"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "78.63"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").Value = "6.486"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("D42").Value = "0.4540"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").Value = "1.997"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").Value = "0.8364"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("D47").Value = "7.601"
$ws.Range("E47").Value = "  +1.66%  "

$ws.Range("D48").Value = "9.861"
$ws.Range("E48").Value = "  +0.83%  "

$ws.Range("D49").Value = "37.67"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D50").Value = "986.48"
$ws.Range("E50").Value = "  +11.46%  "

$ws.Range("E51").Value = "  +4.71%  "
